$wb = $excel.ActiveWorkbook

# Rename sheets (by index, since names are changing)
$wb.Worksheets.Item(1).Name = "GNG_TO-16504777884044206"
$wb.Worksheets.Item(2).Name = "NB_TO-1650477790749374"
$wb.Worksheets.Item(3).Name = "RS_TO-16504777907503755"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504777907974102"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504777908614125"

# Sheet 1 (GNG) - update B2:B5
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650477788368376.csv"
$ws1.Range("B3").Value = "GNG_stims-16504777883874106.csv"
$ws1.Range("B4").Value = "go_stims-16504777883883758.csv"
$ws1.Range("B5").Value = "GNG_stims-16504777884034107.csv"

# Sheet 2 (NB) - update B2:B10
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16504777890973783.csv"
$ws2.Range("B3").Value = "TB-16504777902563782.csv"
$ws2.Range("B4").Value = "ZB-match_1-16504777885853791.csv"
$ws2.Range("B5").Value = "TB-16504777906244104.csv"
$ws2.Range("B6").Value = "OB-16504777889863772.csv"
$ws2.Range("B7").Value = "TB-16504777907303748.csv"
$ws2.Range("B8").Value = "ZB-match_2-1650477788873378.csv"
$ws2.Range("B9").Value = "ZB-match_9-1650477788611381.csv"
$ws2.Range("B10").Value = "OB-16504777896884058.csv"

# Sheet 4 (TOL) - update B2:B7
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504777907653773.csv"
$ws4.Range("B3").Value = "ZM_stims-16504777907523766.csv"
$ws4.Range("B4").Value = "MM_stims-1650477790781413.csv"
$ws4.Range("B5").Value = "ZM_stims-16504777907653773.csv"
$ws4.Range("B6").Value = "MM_stims-16504777907974102.csv"
$ws4.Range("B7").Value = "ZM_stims-16504777907823784.csv"

# Sheet 5 (vSAT) - update B2:B5
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16504777908294127.csv"
$ws5.Range("B3").Value = "SAT_stims-16504777908003764.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504777908454113.csv"
$ws5.Range("B5").Value = "SAT_stims-1650477790813412.csv"
